$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hyperlinks in B2/B3 keep pointing at the same target URLs (deaths /
# positives datasets respectively); only the visible cell text changes to
# the new, more descriptive dataset names pulled from the analyzed source.
$ws.Range("B3").Value = "Positive cases due to COVID-19 - [Ministry of Health - MINSA]"
$ws.Range("B2").Value = "Deaths from COVID-19 - [Ministry of Health - MINSA]"
